$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 186 is the last existing data row; append rows 187 and 188 as copies of it,
# with column A (the date) incremented by one day per new row.
$sourceRow = 186
$baseDate = $ws.Cells.Item($sourceRow, 1).Value2

for ($i = 1; $i -le 2; $i++) {
    $newRow = $sourceRow + $i

    # Copy the whole source row's formatting/styles onto the new row first,
    # so the new cells keep the same style index (s="2" on col A, etc.)
    $srcRange = $ws.Range($ws.Cells.Item($sourceRow, 1), $ws.Cells.Item($sourceRow, 10))
    $dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 10))
    $srcRange.Copy($dstRange)

    # Now set the values: column A gets the incremented date, B..J repeat row 186's values.
    $ws.Cells.Item($newRow, 1).Value2 = $baseDate + $i
    for ($col = 2; $col -le 10; $col++) {
        $ws.Cells.Item($newRow, $col).Value2 = $ws.Cells.Item($sourceRow, $col).Value2
    }
}
